# Add a new feedback bullet point after the "Pleased with the idea..." item.
# The new bullet uses the same ListParagraph / numbered-list (numId 1) formatting
# as the existing feedback bullets, and the "_GoBack" bookmark (which Word places
# at the last edited location) should end up at the end of this new, final
# paragraph - matching what a real edit-in-Word session would produce.

$d = $word.ActiveDocument

$newText = "Thought the pitch was well-written and short, despite their highly critical nature. Agree as long as the website is simple and easy to get around. "

# Temporarily remove the existing _GoBack bookmark; we'll recreate it in the
# right spot once the new paragraph/text exists, rather than letting Word's
# default bookmark-stickiness leave it behind in the old last paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the current last paragraph ("Pleased with the idea...") and insert a
# brand-new paragraph right after its text but before its paragraph mark, so
# the original run (and its rsid attribute) is left completely untouched.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$insertPoint = $d.Range($lastRange.End - 1, $lastRange.End - 1)
$insertPoint.InsertParagraphAfter()

# The newly inserted paragraph automatically inherits the ListParagraph /
# numPr formatting from the paragraph it split off from. Fill in its text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = $newText

# Re-anchor the _GoBack bookmark at the very end of the new last paragraph
# (after the text, before the paragraph mark). A zero-width bookmark Range
# placed at the absolute end of the document body doesn't always round-trip
# correctly, so we bracket the spot with a throwaway character while adding
# the bookmark, then remove that character again.
$newLast = $d.Paragraphs.Last
$targetPos = $newLast.Range.End - 1

$tempRange = $d.Range($targetPos, $targetPos)
$tempRange.InsertAfter("X")

$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tempCharRange = $d.Range($targetPos, $targetPos + 1)
$tempCharRange.Delete()
